$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "298.22") must be
# pre-formatted as Text, otherwise Excel auto-converts the literal into a
# floating point number (losing the exact decimal text / introducing FP noise).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.105.42"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "2.260.15"
$ws.Range("E3").Value = "  -3.45%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "298.22"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "93.97"
$ws.Range("E6").Value = "  -6.69%  "
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -3.78%  "
$ws.Range("D10").Value = "33.01"
$ws.Range("E10").Value = "  -5.74%  "
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "47.79"
$ws.Range("E12").Value = "  -8.40%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "2.609.56"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "15.25"
$ws.Range("E16").Value = "  -3.81%  "
$ws.Range("D17").Value = "2.259.07"
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").Value = "0.778"
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("D19").Value = "42.110.45"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").Value = "0.0₃0888"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").Value = "11.38"
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("E22").Value = "  -3.80%  "
$ws.Range("D23").Value = "66.57"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").Value = "233.13"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "1.92"
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -4.13%  "
$ws.Range("D28").Value = "23.74"
$ws.Range("E28").Value = "  -7.22%  "
$ws.Range("D29").Value = "166.84"
$ws.Range("E29").Value = "  +4.45%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "33.61"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -12.66%  "
$ws.Range("D32").Value = "9.04"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "4.40"
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.0693"
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("D38").Value = "2.80"
$ws.Range("E38").Value = "  -5.65%  "
$ws.Range("D39").Value = "15.99"
$ws.Range("E39").Value = "  -7.92%  "
$ws.Range("D40").Value = "0.0985"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("E42").Value = "  -8.35%  "
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "1.938.45"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("D45").Value = "0.0279"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("E46").Value = "  -7.83%  "
$ws.Range("D47").Value = "9.54"
$ws.Range("E47").Value = "  -7.56%  "
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -5.64%  "
$ws.Range("D49").Value = "2.81"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "2.483.94"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("D51").Value = "52.28"
$ws.Range("E51").Value = "  -7.13%  "
